# Add a new "toi_source" attribute row to the ncp/gop ColumnHeaders sheet.
# This inserts a new row 6 (pushing gop/ncp/ncp_per_gop down to rows 7-9)
# and fills only the attributeName / attributeDefinition / class columns
# for the new attribute, matching the commit "add ncp toi_source, remove
# temp/sal from output".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 6 (gop).
$ws.Rows("6:6").Insert()

# The inserted row inherits formatting from the row above (B5 has a
# wrap-text style); strip that so the new row matches a plain/default row.
$ws.Rows("6:6").ClearFormats()

# Populate the new attribute row.
$ws.Range("A6").Value = "toi_source"
$ws.Range("B6").Value = "Bottle sample from niskin or underway"
$ws.Range("C6").Value = "categorical"

# Match the author's resulting selection state.
[void]$ws.Range("A6:H6").Select()
